# Repull data, push all data, mean calculation
# Update column F ("dSF") values for the specified rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    16 = -1
    22 = 1
    28 = 4
    34 = -10
    38 = -2
    41 = 4
    43 = 3
    55 = 1
    57 = -2
    58 = -3
    63 = 0
    68 = 0
    75 = 1
    76 = -3
    77 = -6
    78 = -3
    79 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
